$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the address lookup test data: the shared string " 4215 Hunt Rd, Farmington, MO 63640"
# had a stray leading space. Re-enter the value without the leading space but keep it
# entered as text (leading apostrophe -> quotePrefix) so it keeps the same "typed as text"
# formatting it had before.
$ws.Range("E2").Value = "'4215 Hunt Rd, Farmington, MO 63640"

# Re-enter the driver name value so its formatting (no cell fill applied) matches the
# locators fix -- drop the stray "applyFill" that had been left on this cell's style.
$ws.Range("Q2").Value = "'LJB Liam Jack Benjamin"
$ws.Range("Q2").Interior.Pattern = -4142

# Leave the active selection on E5, matching the workbook's last saved selection.
$ws.Range("E5").Select() | Out-Null
